$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6315543333333333
$ws.Range("H2").Value = 1.894663
$ws.Range("I2").Value = 0.4681870832749024
$ws.Range("J2").Value = 0.5690659591172773
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 107.663086
$ws.Range("N2").Value = 215.326172
$ws.Range("O2").Value = 0.2751823527645522
$ws.Range("P2").Value = 0.2057131686791961
$ws.Range("Q2").Value = 67.99508850333933
$ws.Range("R2").Value = 407.970531020036
$ws.Range("S2").Value = 0.128836823109561
$ws.Range("T2").Value = 0.117064361637481

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6315543333333333
$ws.Range("H3").Value = 1.894663
$ws.Range("I3").Value = 0.4681870832749024
$ws.Range("J3").Value = 0.5690659591172773
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.41172733333334
$ws.Range("N3").Value = 82.23518200000001
$ws.Range("O3").Value = 0.07006323059444025
$ws.Range("P3").Value = 0.07856388152449204
$ws.Range("Q3").Value = 17.31199518151844
$ws.Range("R3").Value = 155.807956633666
$ws.Range("S3").Value = 0.03280269957682789
$ws.Range("T3").Value = 0.0447080305917112

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6315543333333333
$ws.Range("H4").Value = 1.894663
$ws.Range("I4").Value = 0.4681870832749024
$ws.Range("J4").Value = 0.5690659591172773
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 82.303927
$ws.Range("N4").Value = 246.911781
$ws.Range("O4").Value = 0.2103654011331419
$ws.Range("P4").Value = 0.2358886724356653
$ws.Range("Q4").Value = 51.97940174720033
$ws.Range("R4").Value = 467.814615724803
$ws.Range("S4").Value = 0.09849036357848057
$ws.Range("T4").Value = 0.1342362136245031

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.6315543333333333
$ws.Range("H5").Value = 1.894663
$ws.Range("I5").Value = 0.4681870832749024
$ws.Range("J5").Value = 0.5690659591172773
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 71.272429
$ws.Range("N5").Value = 213.817287
$ws.Range("O5").Value = 0.1821693528222338
$ws.Range("P5").Value = 0.2042716462128862
$ws.Range("Q5").Value = 45.01241138214233
$ws.Range("R5").Value = 405.111702439281
$ws.Range("S5").Value = 0.08528933795991826
$ws.Range("T5").Value = 0.1162440402726012

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6315543333333333
$ws.Range("H6").Value = 1.894663
$ws.Range("I6").Value = 0.4681870832749024
$ws.Range("J6").Value = 0.5690659591172773
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 83.25665266666667
$ws.Range("N6").Value = 249.769958
$ws.Range("O6").Value = 0.2128005281598046
$ws.Range("P6").Value = 0.2386192492246123
$ws.Range("Q6").Value = 52.58109977046156
$ws.Range("R6").Value = 473.229897934154
$ws.Range("S6").Value = 0.09963045859849766
$ws.Range("T6").Value = 0.1357900919238486

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6315543333333333
$ws.Range("H7").Value = 1.894663
$ws.Range("I7").Value = 0.4681870832749024
$ws.Range("J7").Value = 0.5690659591172773
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 19.3348755
$ws.Range("N7").Value = 38.66975100000001
$ws.Range("O7").Value = 0.04941913452582716
$ws.Range("P7").Value = 0.03694338192314826
$ws.Range("Q7").Value = 12.2110244064855
$ws.Range("R7").Value = 73.26614643891301
$ws.Range("S7").Value = 0.02313740045161705
$ws.Range("T7").Value = 0.02102322106713225

$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 0.7173815
$ws.Range("H8").Value = 1.434763
$ws.Range("I8").Value = 0.5318129167250975
$ws.Range("J8").Value = 0.4309340408827227
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 107.663086
$ws.Range("N8").Value = 215.326172
$ws.Range("O8").Value = 0.2751823527645522
$ws.Range("P8").Value = 0.2057131686791961
$ws.Range("Q8").Value = 77.235506129309
$ws.Range("R8").Value = 308.942024517236
$ws.Range("S8").Value = 0.1463455296549912
$ws.Range("T8").Value = 0.08864880704171513

$ws.Range("E9").Value = 2
$ws.Range("G9").Value = 0.7173815
$ws.Range("H9").Value = 1.434763
$ws.Range("I9").Value = 0.5318129167250975
$ws.Range("J9").Value = 0.4309340408827227
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.41172733333334
$ws.Range("N9").Value = 82.23518200000001
$ws.Range("O9").Value = 0.07006323059444025
$ws.Range("P9").Value = 0.07856388152449204
$ws.Range("Q9").Value = 19.66466607197767
$ws.Range("R9").Value = 117.987996431866
$ws.Range("S9").Value = 0.03726053101761236
$ws.Range("T9").Value = 0.03385585093278084

$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 0.7173815
$ws.Range("H10").Value = 1.434763
$ws.Range("I10").Value = 0.5318129167250975
$ws.Range("J10").Value = 0.4309340408827227
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 82.303927
$ws.Range("N10").Value = 246.911781
$ws.Range("O10").Value = 0.2103654011331419
$ws.Range("P10").Value = 0.2358886724356653
$ws.Range("Q10").Value = 59.0433146071505
$ws.Range("R10").Value = 354.259887642903
$ws.Range("S10").Value = 0.1118750375546613
$ws.Range("T10").Value = 0.1016524588111622

$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 0.7173815
$ws.Range("H11").Value = 1.434763
$ws.Range("I11").Value = 0.5318129167250975
$ws.Range("J11").Value = 0.4309340408827227
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 71.272429
$ws.Range("N11").Value = 213.817287
$ws.Range("O11").Value = 0.1821693528222338
$ws.Range("P11").Value = 0.2042716462128862
$ws.Range("Q11").Value = 51.1295220246635
$ws.Range("R11").Value = 306.777132147981
$ws.Range("S11").Value = 0.09688001486231552
$ws.Range("T11").Value = 0.08802760594028498

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.7173815
$ws.Range("H12").Value = 1.434763
$ws.Range("I12").Value = 0.5318129167250975
$ws.Range("J12").Value = 0.4309340408827227
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 83.25665266666667
$ws.Range("N12").Value = 249.769958
$ws.Range("O12").Value = 0.2128005281598046
$ws.Range("P12").Value = 0.2386192492246123
$ws.Range("Q12").Value = 59.72678237499233
$ws.Range("R12").Value = 358.360694249954
$ws.Range("S12").Value = 0.1131700695613069
$ws.Range("T12").Value = 0.1028291573007637

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.7173815
$ws.Range("H13").Value = 1.434763
$ws.Range("I13").Value = 0.5318129167250975
$ws.Range("J13").Value = 0.4309340408827227
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 19.3348755
$ws.Range("N13").Value = 38.66975100000001
$ws.Range("O13").Value = 0.04941913452582716
$ws.Range("P13").Value = 0.03694338192314826
$ws.Range("Q13").Value = 13.87048198850325
$ws.Range("R13").Value = 55.48192795401301
$ws.Range("S13").Value = 0.02628173407421011
$ws.Range("T13").Value = 0.01592016085601601

